$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on price cells whose new values would otherwise
# be auto-coerced to numbers by Excel (losing the trailing-zero / format).
$textCells = @("D5", "D6", "D10", "D15", "D16", "D18", "D19", "D22", "D23", "D24", "D29", "D37", "D44", "D46", "D47")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.924.36"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.665.38"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "215.62"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "0.532"
$ws.Range("E6").Value = "  +4.70%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").Value = "20.30"
$ws.Range("E10").Value = "  +3.18%  "
$ws.Range("E11").Value = "  +3.86%  "
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "1.662.61"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.525"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "66.22"
$ws.Range("E16").Value = "  +2.11%  "
$ws.Range("D17").Value = "26.918.01"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "233.62"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").Value = "8.01"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "0.0₃0732"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").Value = "2.22"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "9.14"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "1.460.07"
$ws.Range("E33").Value = "  -4.40%  "
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "0.578"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  +6.37%  "
$ws.Range("D44").Value = "65.93"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "1.809.05"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "90.47"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  +4.06%  "
